# Commit EPG and PIP Test cases for Release
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) screenTitles sheet ("sheet4.xml"): add column C (name_fr) for the first
#    4 rows. Values are added in column order (C1,C2,C3) so new shared
#    strings land at the expected indices; C4 re-uses the existing "epg"
#    string already used by B4.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("screenTitles")

$ws4.Range("C1").Value = "name_fr"
$ws4.Range("C2").Value = "paramètres"
$ws4.Range("C3").Value = "système"
$ws4.Range("C4").Value = "epg"

# match header formatting (bold / fill / alignment) used by A1:B1
$ws4.Range("B1").Copy()
$ws4.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# approximate the bestFit column width recorded for the new column
$ws4.Columns("C").ColumnWidth = 10.333333333

$ws4.Range("G4").Select()

# ---------------------------------------------------------------------------
# 2) Add the new "PIPScreen" worksheet at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add($null, $lastSheet)
$ws8.Name = "PIPScreen"

# Values are populated in this particular order (bottom-right to top-left)
# so that the shared-string table ends up with the same ordering as the
# authored workbook.
$ws8.Range("C3").Value = "Right"
$ws8.Range("B3").Value = "links"
$ws8.Range("A3").Value = "PIPLink"
$ws8.Range("C2").Value = "Left"
$ws8.Range("B2").Value = "standaard"
$ws8.Range("A2").Value = "Default"
$ws8.Range("C1").Value = "pip_position"
$ws8.Range("B1").Value = "name_nl"
$ws8.Range("A1").Value = "objectID"

# match header formatting on row 1
$ws4.Range("A1:B1").Copy()
$ws8.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws8.Range("H18").Select()

# ---------------------------------------------------------------------------
# 3) Re-select the screenTitles sheet as the active tab (it was the active
#    tab in the saved workbook).
# ---------------------------------------------------------------------------
$ws4.Activate()
